$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: add "Comp" / "Compare" in A14:B14, matching the style used in row 10 (A10:B10) ---
$ws.Range("A10:B10").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = "Comp"
$ws.Range("B14").Value = "Compare"

# --- Row 15: add "ALL" / "all avaiable waves in ABCD 5.0" in A15:B15, matching the style used in row 13 (A13:B13) ---
$ws.Range("A13:B13").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15").Value = "ALL"
$ws.Range("B15").Value = "all avaiable waves in ABCD 5.0"

# Row 15 grew taller to fit the wrapped description text
$ws.Rows("15").RowHeight = 34.5

$excel.CutCopyMode = 0

# Update the active selection to match the author's final cursor position
$ws.Range("G13").Select() | Out-Null
